$d = $word.ActiveDocument

# Locate the email address run in the contact-info paragraph so we can
# anchor the deletion right after it (keeping "Email: <link>" intact).
$emailRange = $d.Content
$found = $emailRange.Find.Execute("kristyntaniguchi@yahoo.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $emailEnd = $emailRange.End

    # The paragraph's end position (Range.End includes the trailing
    # paragraph mark), so stop one character short of it.
    $para = $emailRange.Paragraphs(1)
    $paraEnd = $para.Range.End

    if ($paraEnd - 1 -gt $emailEnd) {
        $trailRange = $d.Range($emailEnd, $paraEnd - 1)
        $trailRange.Delete()
    }
}
